{"js": "// The document begins with a title paragraph \"Featured Images\" followed by\n// TWO identical empty, centered, Helvetica 48pt, underlined paragraphs\n// before the body text begins. The edit removes one of these duplicate\n// empty paragraphs, leaving a single empty paragraph between the title and\n// the body text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the title paragraph (\"Featured Images\").\nlet titleIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"Featured Images\") {\n    titleIndex = i;\n    break;\n  }\n}\n\nlet target = null;\n\nif (titleIndex !== -1) {\n  // Prefer the paragraph immediately after the title if it is empty (the\n  // first of the two duplicate empty paragraphs introduced by the diff).\n  if (titleIndex + 1 < items.length && items[titleIndex + 1].text === \"\") {\n    target = items[titleIndex + 1];\n  }\n}\n\nif (!target) {\n  // Fallback: find the first pair of consecutive empty paragraphs and\n  // remove the first one of the pair.\n  for (let i = 0; i + 1 < items.length; i++) {\n    if (items[i].text === \"\" && items[i + 1].text === \"\") {\n      target = items[i];\n      break;\n    }\n  }\n}\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# The document begins with a title paragraph \"Featured Images\" followed by\n# TWO identical empty, centered, Helvetica 48pt, underlined paragraphs\n# before the body text begins. The edit removes one of these duplicate\n# empty paragraphs, leaving a single empty paragraph between the title and\n# the body text.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaText($para) {\n    return $para.Range.Text.TrimEnd(\"`r\")\n}\n\n$count = $d.Paragraphs.Count\n\n# Locate the title paragraph (\"Featured Images\").\n$titleIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ((Get-ParaText $d.Paragraphs.Item($i)) -eq \"Featured Images\") {\n        $titleIndex = $i\n        break\n    }\n}\n\n$targetIndex = -1\n\nif ($titleIndex -ne -1) {\n    $nextIndex = $titleIndex + 1\n    if ($nextIndex -le $count) {\n        if ((Get-ParaText $d.Paragraphs.Item($nextIndex)) -eq \"\") {\n            $targetIndex = $nextIndex\n        }\n    }\n}\n\nif ($targetIndex -eq -1) {\n    # Fallback: find the first pair of consecutive empty paragraphs and\n    # remove the first one of the pair.\n    for ($i = 1; $i -lt $count; $i++) {\n        if ((Get-ParaText $d.Paragraphs.Item($i)) -eq \"\" -and (Get-ParaText $d.Paragraphs.Item($i + 1)) -eq \"\") {\n            $targetIndex = $i\n            break\n        }\n    }\n}\n\nif ($targetIndex -ne -1) {\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n}\n"}
